# Update Argent (Silver) prices in the "Solar_Prices" workbook.
# Each target sheet stores its latest price as text in cell B5; a leading
# apostrophe forces Excel to keep the value as literal text instead of
# re-interpreting it as a number (important for comma thousand-separators
# and decimal values).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cell Topcon 183mm")
$ws.Range("B5").Value = "'0.295"

$ws = $wb.Worksheets.Item("Silver Rear_side")
$ws.Range("B5").Value = "'5,282"

$ws = $wb.Worksheets.Item("Silver Busbar front-side")
$ws.Range("B5").Value = "'7,907"

$ws = $wb.Worksheets.Item("Silver finger front-side")
$ws.Range("B5").Value = "'7,957"

$ws = $wb.Worksheets.Item("USD_CNY")
$ws.Range("B5").Value = "'7.2617"
